$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "SRPbVT"
# ------------------------------------------------------------------
$srpbvt = $wb.Worksheets.Item("SRPbVT")
$srpbvt.Activate()

# New leading "(Boolean)" header above the pollutant codes
$srpbvt.Range("A1").Value = "(Boolean)"

# Aircraft / Rail no longer flagged for SOx (column G)
$srpbvt.Range("G5").Value = 0
$srpbvt.Range("G6").Value = 0

[void]$srpbvt.Range("G7").Select()

# ------------------------------------------------------------------
# Sheet "About"
# ------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Activate()

# Fix/introduce wording tweak in the VOC note (named -> names)
$about.Range("A36").Value = "Although VOCs are not specifically named as a criteria pollutant, ozone is names, and ozone is formed from"

# Turn the EPA marine-engine source URL in B26 into a real hyperlink
$about.Hyperlinks.Add($about.Range("B26"), "https://nepis.epa.gov/Exe/ZyPDF.cgi/P1005ZAD.PDF?Dockey=P1005ZAD.PDF")

# Restore the current selection/scroll position that was saved with the file
[void]$about.Range("B27").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
